# Applies the "Add files via upload" edit:
#  - Updates the Last Update date (A2) to 3/27/2019
#  - Inserts 6 new rows of meeting-attendance data in the
#    "Regularly Scheduled Meetings" block (new rows 48-53)
#  - Adds a new "Ad hoc Meetings" entry (new row 64, "Recevie testing")
#  - Moves the selection / view to match the saved file

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "Last Update" date ---------------------------------------
$ws.Range("A2").Value = (Get-Date -Year 2019 -Month 3 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# --- Insert 6 new rows before the old row 48 --------------------------
$ws.Rows("48:53").Insert()

# Row 48: Team meeting, Wednesday, March 11, 2019
$ws.Range("A48").Value = "Team"
$ws.Range("B48").Value = "Wednesday, March 11, 2019"
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 1

# Row 49: Faculty Adv meeting, Wednesday, March 13, 2019
$ws.Range("A49").Value = "Faculty Adv"
$ws.Range("B49").Value = "Wednesday, March 13, 2019"
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 1

# Row 50: Team meeting, Wednesday, March 13, 2019
$ws.Range("A50").Value = "Team"
$ws.Range("B50").Value = "Wednesday, March 13, 2019"
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 1

# Row 51: Team meeting, Wednesday, March 25, 2019
$ws.Range("A51").Value = "Team"
$ws.Range("B51").Value = "Wednesday, March 25, 2019"
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 1

# Row 52: Team meeting, Wednesday, March 27, 2019
$ws.Range("A52").Value = "Team"
$ws.Range("B52").Value = "Wednesday, March 27, 2019"
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 1

# Row 53: Faculty Adv meeting, Wednesday, March 27, 2019
$ws.Range("A53").Value = "Faculty Adv"
$ws.Range("B53").Value = "Wednesday, March 27, 2019"
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 1
$ws.Range("E53").Value = 1

# Apply formatting (style) consistent with neighbouring rows for B:E
$ws.Range("B48:B53").Style = $ws.Range("B47").Style
$ws.Range("C48:E53").Style = $ws.Range("C47:E47").Style

# --- Fill in the new "Audio Pass through" result (new row 63) ---------
$ws.Range("C63").Value = 1
$ws.Range("D63").Value = 1
$ws.Range("E63").Value = 1

# --- Add the new "Recevie testing" ad-hoc meeting (new row 64) --------
$ws.Range("A64").Value = "Recevie testing"
$ws.Range("B64").Value = "Sunday, March 30, 2019"
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 1
$ws.Range("E64").Value = 1
$ws.Range("A64:E64").Style = $ws.Range("A63:E63").Style

# --- Update view / selection to match the saved state -----------------
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("I57").Select()
